# Insert a new weekly data row into the "Arveja Verde" sheet.
# The new row is inserted at row 58 (pushing the existing rows 58-79 down
# to 59-80), and its values are filled in below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 58 (shifts rows 58..79 down to 59..80)
$ws.Rows.Item(58).Insert()

# Fill in the values for the newly inserted row 58
$ws.Cells.Item(58, 1).Value = 4
$ws.Cells.Item(58, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(58, 3).Value = "Los Lagos"
$ws.Cells.Item(58, 4).Value = 44510
$ws.Cells.Item(58, 5).Value = 10
$ws.Cells.Item(58, 6).Value = 100112022
$ws.Cells.Item(58, 7).Value = "Arveja Verde"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 90
$ws.Cells.Item(58, 11).Value = 20000
$ws.Cells.Item(58, 12).Value = 20000
$ws.Cells.Item(58, 13).Value = 20000
$ws.Cells.Item(58, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(58, 15).Value = "Región Metropolitana"
$ws.Cells.Item(58, 16).Value = 800
$ws.Cells.Item(58, 17).Value = 25
$ws.Cells.Item(58, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date style as the other rows (s="2")
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
